$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume(1h) (E) columns keep their original text formatting
# so that numeric-looking strings are not auto-converted to numbers/dates.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "59.186.59"
$ws.Range("E2").Value = "  -1.38%  "
$ws.Range("D3").Value = "2.739.13"
$ws.Range("E3").Value = "  -5.30%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "509.46"
$ws.Range("E5").Value = "  -2.73%  "
$ws.Range("D6").Value = "141.97"
$ws.Range("E6").Value = "  +0.34%  "
$ws.Range("E7").Value = "  -0.42%  "
$ws.Range("D8").Value = "0.535"
$ws.Range("E8").Value = "  -2.33%  "
$ws.Range("D9").Value = "2.749.23"
$ws.Range("E9").Value = "  -5.25%  "
$ws.Range("D10").Value = "6.14"
$ws.Range("E10").Value = "  +5.08%  "
$ws.Range("D11").Value = "0.105"
$ws.Range("E11").Value = "  -1.04%  "
$ws.Range("D12").Value = "0.351"
$ws.Range("E12").Value = "  -0.65%  "
$ws.Range("E13").Value = "  +1.59%  "
$ws.Range("D14").Value = "3.212.88"
$ws.Range("E14").Value = "  -5.76%  "
$ws.Range("D15").Value = "59.021.62"
$ws.Range("E15").Value = "  -1.94%  "
$ws.Range("D16").Value = "21.94"
$ws.Range("E16").Value = "  -2.55%  "
$ws.Range("D17").Value = "0.0000137"
$ws.Range("E17").Value = "  -1.15%  "
$ws.Range("D18").Value = "2.727.58"
$ws.Range("E18").Value = "  -6.05%  "
$ws.Range("D19").Value = "4.77"
$ws.Range("E19").Value = "  -2.70%  "
$ws.Range("D20").Value = "11.06"
$ws.Range("E20").Value = "  -3.30%  "
$ws.Range("D21").Value = "348.04"
$ws.Range("E21").Value = "  -2.81%  "
$ws.Range("D22").Value = "6.29"
$ws.Range("E22").Value = "  -3.50%  "
$ws.Range("E23").Value = "  +0.22%  "
$ws.Range("D24").Value = "5.62"
$ws.Range("E24").Value = "  -0.26%  "
$ws.Range("D25").Value = "63.15"
$ws.Range("E25").Value = "  +0.13%  "
$ws.Range("D26").Value = "0.427"
$ws.Range("E26").Value = "  -4.39%  "
$ws.Range("D27").Value = "0.174"
$ws.Range("E27").Value = "  +0.89%  "
$ws.Range("D28").Value = "1.02"
$ws.Range("E28").Value = "  +1.73%  "
$ws.Range("D29").Value = "0.0₃0843"
$ws.Range("E29").Value = "  +0.34%  "
$ws.Range("D30").Value = "7.53"
$ws.Range("E30").Value = "  -2.35%  "
$ws.Range("E31").Value = "  -0.16%  "
$ws.Range("D32").Value = "1.62"
$ws.Range("E32").Value = "  -1.62%  "
$ws.Range("D33").Value = "19.22"
$ws.Range("E33").Value = "  -0.66%  "
$ws.Range("D34").Value = "149.67"
$ws.Range("E34").Value = "  -0.32%  "
$ws.Range("D35").Value = "4.25"
$ws.Range("E35").Value = "  -1.10%  "
$ws.Range("D36").Value = "5.42"
$ws.Range("E36").Value = "  -1.42%  "
$ws.Range("D37").Value = "0.955"
$ws.Range("E37").Value = "  -2.41%  "
$ws.Range("D38").Value = "1.15"
$ws.Range("E38").Value = "  -2.83%  "
$ws.Range("D39").Value = "36.24"
$ws.Range("E39").Value = "  -3.78%  "
$ws.Range("D40").Value = "1.40"
$ws.Range("E40").Value = "  -3.82%  "
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").Value = "3.56"
$ws.Range("E41").Value = "  -2.28%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "2.193.96"
$ws.Range("E42").Value = "  -5.85%  "
$ws.Range("D43").Value = "0.0562"
$ws.Range("E43").Value = "  -0.62%  "
$ws.Range("D44").Value = "0.994"
$ws.Range("E44").Value = "  -0.64%  "
$ws.Range("D45").Value = "0.605"
$ws.Range("E45").Value = "  -5.71%  "
$ws.Range("D46").Value = "19.21"
$ws.Range("E46").Value = "  -7.07%  "
$ws.Range("D47").Value = "4.83"
$ws.Range("E47").Value = "  +0.58%  "
$ws.Range("E48").Value = "  +0.22%  "
$ws.Range("D49").Value = "0.0231"
$ws.Range("E49").Value = "  -0.04%  "
$ws.Range("D50").Value = "0.0888"
$ws.Range("E50").Value = "  -3.79%  "
$ws.Range("D51").Value = "18.16"
$ws.Range("E51").Value = "  +0.74%  "
